$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 to 7 (Patient-075, Patient-080, Patient-088, Patient-089)
$ws.Range("A4:G7").Delete()

# Update remaining two data rows with new patient names and values
$ws.Range("A2").Value = "Patient-093"
$ws.Range("B2").Value = 86
$ws.Range("C2").Value = 990
$ws.Range("D2").Value = 317
$ws.Range("E2").Value = 639
$ws.Range("F2").Value = 858
$ws.Range("G2").Value = 1111

$ws.Range("A3").Value = "Patient-094"
$ws.Range("B3").Value = 71
$ws.Range("C3").Value = 430
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 122
$ws.Range("F3").Value = 283
$ws.Range("G3").Value = 779

$wb.Save()
